# Generate Report for Handoff
# Updates the localization-status report after a handoff round:
#   - Overview sheet: bump "Latest HO Xliff Generate Date" for the
#     5691be4f-728a-4911-9f5b-eaf70c7f9122.md row family (rows 8-14,
#     skipping row 12 / afa6aa99...) from 06:29:21 -> 06:29:51
#   - zh-cn sheet: same rows get Priority "ht" and their
#     "Latest Handoff Datetime" bumped from 06:29:13 -> 06:29:43
#   - de-de sheet: same rows get Priority "ht" (its "Latest Handoff
#     Datetime" column shares the same string as the Overview sheet's
#     date, so it updates together with the Overview edit above)

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = @(8, 9, 10, 11, 13, 14)

foreach ($r in $rows) {
    # Overview!G<r> : Latest HO Xliff Generate Date
    $overview.Range("G$r").Value = "2016-09-07 06:29:51"

    # zh-cn!E<r> : Priority
    $zhcn.Range("E$r").Value = "ht"
    # zh-cn!H<r> : Latest Handoff Datetime
    $zhcn.Range("H$r").Value = "2016-09-07 06:29:43"

    # de-de!E<r> : Priority
    $dede.Range("E$r").Value = "ht"
    # de-de!H<r> : Latest Handoff Datetime (shares the same underlying
    # text as Overview!G -- keep both in sync)
    $dede.Range("H$r").Value = "2016-09-07 06:29:51"
}
